$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 = "season" attribute row; column F = "Values (range, description)".
# Previously blank -> now documents the actual season values used.
$ws.Range("F13").Value = "fall, spring, summer"

# Row 14 = "estuary" attribute row; column F = "Values (range, description)".
# Extend the existing value list to mention the additional "ere" estuary code.
$ws.Range("F14").Value = "lqre (Little Qualicum River Estuary), nre (Nanaimo River Estuary), ere (Englishman River Estuary, not included in analysis)"

# Reflect the cell that was selected / in view when the author last saved the file.
$ws.Range("F15").Select()
